$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.68"
$ws.Range("E2").Value = "'-4.68%"
$ws.Range("D3").Value = "'39.73"
$ws.Range("E3").Value = "'-7.87%"
$ws.Range("D4").Value = "'5.120"
$ws.Range("E4").Value = "'-1.78%"
$ws.Range("D5").Value = "'0.07701"
$ws.Range("E5").Value = "'-6.00%"
$ws.Range("D6").Value = "'4.242"
$ws.Range("E6").Value = "'-1.87%"
$ws.Range("D7").Value = "'1.622"
$ws.Range("E7").Value = "'-11.87%"
$ws.Range("D8").Value = "'0.8800"
$ws.Range("E8").Value = "'-5.85%"
$ws.Range("D9").Value = "'0.1002"
$ws.Range("E9").Value = "'-10.06%"
$ws.Range("D10").Value = "'0.1746"
$ws.Range("E10").Value = "'-6.64%"
$ws.Range("D11").Value = "'0.08912"
$ws.Range("E11").Value = "'-6.04%"
$ws.Range("D12").Value = "'0.04385"
$ws.Range("E12").Value = "'-5.14%"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.001255"
$ws.Range("E14").Value = "'-2.79%"
$ws.Range("D15").Value = "'0.005902"
$ws.Range("E15").Value = "'2.81%"
$ws.Range("D16").Value = "'3.352"
$ws.Range("E16").Value = "'-0.19%"
$ws.Range("D17").Value = "'2.437"
$ws.Range("E17").Value = "'-2.95%"
$ws.Range("E18").Value = "'-0.46%"
$ws.Range("D19").Value = "'7.008"
$ws.Range("E19").Value = "'-5.58%"
$ws.Range("D20").Value = "'0.1340"
$ws.Range("E20").Value = "'-3.51%"
$ws.Range("D21").Value = "'0.3003"
$ws.Range("E21").Value = "'14.47%"
$ws.Range("D22").Value = "'0.04157"
$ws.Range("D23").Value = "'0.001200"
$ws.Range("E23").Value = "'-3.83%"
$ws.Range("D24").Value = "'0.004067"
$ws.Range("E24").Value = "'-5.52%"
$ws.Range("D25").Value = "'0.0001221"
$ws.Range("E25").Value = "'11.00%"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("D38").Value = "'0.02336"
$ws.Range("E38").Value = "'-13.97%"
$ws.Range("E39").Value = "'-7.39%"
$ws.Range("D40").Value = "'0.007923"
$ws.Range("E40").Value = "'-0.60%"
$ws.Range("D41").Value = "'0.1321"
$ws.Range("E41").Value = "'-5.34%"
$ws.Range("D42").Value = "'0.006330"
$ws.Range("E42").Value = "'-3.26%"
$ws.Range("D43").Value = "'0.001951"
$ws.Range("E43").Value = "'-6.77%"
$ws.Range("D44").Value = "'0.008495"
$ws.Range("E44").Value = "'13.78%"
$ws.Range("D45").Value = "'0.3056"
$ws.Range("E45").Value = "'-4.62%"
$ws.Range("D46").Value = "'0.00006510"
$ws.Range("E46").Value = "'-6.70%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.007006"
$ws.Range("E48").Value = "'98.49%"
$ws.Range("D49").Value = "'0.004636"
$ws.Range("E49").Value = "'33.83%"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E51").Value = "'0.06%"
